# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the per-sheet Leve profit tables (currentAveragePrice*, LevePrice*,
# LeveProfit* columns) across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 719.76
$ws.Range("I15").Value = 719.76
$ws.Range("K15").Value = 2159.28
$ws.Range("M15").Value = -1990.28

# Row 76
$ws.Range("H76").Value = 3846.7307
$ws.Range("I76").Value = 3800.75
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3800.75
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3485.75
$ws.Range("N76").Value = -4630

# Row 79
$ws.Range("H79").Value = 3846.7307
$ws.Range("I79").Value = 3800.75
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3800.75
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2708.75
$ws.Range("N79").Value = -6184

# Row 82
$ws.Range("H82").Value = 5894.0835
$ws.Range("I82").Value = 3460.4
$ws.Range("J82").Value = 7632.4287
$ws.Range("K82").Value = 10381.2
$ws.Range("L82").Value = 22897.2861
$ws.Range("M82").Value = -9975.200000000001
$ws.Range("N82").Value = -23709.2861

# Row 85
$ws.Range("H85").Value = 5894.0835
$ws.Range("I85").Value = 3460.4
$ws.Range("J85").Value = 7632.4287
$ws.Range("K85").Value = 10381.2
$ws.Range("L85").Value = 22897.2861
$ws.Range("M85").Value = -8977.200000000001
$ws.Range("N85").Value = -25705.2861

# Row 86
$ws.Range("H86").Value = 3655.7144
$ws.Range("I86").Value = 2918
$ws.Range("J86").Value = 5500
$ws.Range("K86").Value = 2918
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = -1795
$ws.Range("N86").Value = -7746

# Row 88
$ws.Range("H88").Value = 2139.6155
$ws.Range("I88").Value = 2301.5
$ws.Range("J88").Value = 1600
$ws.Range("K88").Value = 2301.5
$ws.Range("L88").Value = 1600
$ws.Range("M88").Value = -1895.5
$ws.Range("N88").Value = -2412

# Row 89
$ws.Range("H89").Value = 3655.7144
$ws.Range("I89").Value = 2918
$ws.Range("J89").Value = 5500
$ws.Range("K89").Value = 14590
$ws.Range("L89").Value = 27500
$ws.Range("M89").Value = -8974
$ws.Range("N89").Value = -38732

# Row 91
$ws.Range("H91").Value = 2139.6155
$ws.Range("I91").Value = 2301.5
$ws.Range("J91").Value = 1600
$ws.Range("K91").Value = 2301.5
$ws.Range("L91").Value = 1600
$ws.Range("M91").Value = -897.5
$ws.Range("N91").Value = -4408

# Row 93
$ws.Range("H93").Value = 39601
$ws.Range("J93").Value = 39601
$ws.Range("L93").Value = 39601
$ws.Range("N93").Value = -44593

# Row 98
$ws.Range("H98").Value = 943.6842
$ws.Range("I98").Value = 683.125
$ws.Range("J98").Value = 2333.3333
$ws.Range("K98").Value = 683.125
$ws.Range("L98").Value = 2333.3333
$ws.Range("M98").Value = 814.875
$ws.Range("N98").Value = -5329.3333

# Row 122
$ws.Range("H122").Value = 943.6842
$ws.Range("I122").Value = 683.125
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 2049.375
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = 400.625
$ws.Range("N122").Value = -11899.9999

# Row 129
$ws.Range("H129").Value = 1416.2192
$ws.Range("I129").Value = 461.16666
$ws.Range("J129").Value = 1604.0984
$ws.Range("K129").Value = 1383.49998
$ws.Range("L129").Value = 4812.2952
$ws.Range("M129").Value = 3616.50002
$ws.Range("N129").Value = -14812.2952

# Row 135
$ws.Range("H135").Value = 1376.9722
$ws.Range("I135").Value = 1406.3667
$ws.Range("J135").Value = 1230
$ws.Range("K135").Value = 12657.3003
$ws.Range("L135").Value = 11070
$ws.Range("M135").Value = -10122.3003
$ws.Range("N135").Value = -16140

# Row 138
$ws.Range("H138").Value = 5052984
$ws.Range("I138").Value = 9525378
$ws.Range("J138").Value = 3506.6453
$ws.Range("K138").Value = 28576134
$ws.Range("L138").Value = 10519.9359
$ws.Range("M138").Value = -28570994
$ws.Range("N138").Value = -20799.9359

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1992.0294
$ws.Range("I74").Value = 1291.6666
$ws.Range("J74").Value = 3123.3845
$ws.Range("K74").Value = 1291.6666
$ws.Range("L74").Value = 3123.3845
$ws.Range("M74").Value = -417.6666
$ws.Range("N74").Value = -4871.3845

# Row 77
$ws.Range("H77").Value = 1992.0294
$ws.Range("I77").Value = 1291.6666
$ws.Range("J77").Value = 3123.3845
$ws.Range("K77").Value = 6458.333000000001
$ws.Range("L77").Value = 15616.9225
$ws.Range("M77").Value = -2090.333000000001
$ws.Range("N77").Value = -24352.9225

# Row 102
$ws.Range("H102").Value = 2228
$ws.Range("I102").Value = 1819.6
$ws.Range("J102").Value = 3249
$ws.Range("K102").Value = 1819.6
$ws.Range("L102").Value = 3249
$ws.Range("M102").Value = -197.5999999999999
$ws.Range("N102").Value = -6493

# Row 132
$ws.Range("H132").Value = 4238902.5
$ws.Range("I132").Value = 5557002.5
$ws.Range("J132").Value = 2152.7144
$ws.Range("K132").Value = 16671007.5
$ws.Range("L132").Value = 6458.1432
$ws.Range("M132").Value = -16668477.5
$ws.Range("N132").Value = -11518.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2883.3635
$ws.Range("I20").Value = 3141.6667
$ws.Range("J20").Value = 2573.4
$ws.Range("K20").Value = 3141.6667
$ws.Range("L20").Value = 2573.4
$ws.Range("M20").Value = -2894.6667
$ws.Range("N20").Value = -3067.4

# Row 86
$ws.Range("H86").Value = 18520956
$ws.Range("I86").Value = 2242.9524
$ws.Range("J86").Value = 83336450
$ws.Range("K86").Value = 2242.9524
$ws.Range("L86").Value = 83336450
$ws.Range("M86").Value = -1119.9524
$ws.Range("N86").Value = -83338696

# Row 89
$ws.Range("H89").Value = 18520956
$ws.Range("I89").Value = 2242.9524
$ws.Range("J89").Value = 83336450
$ws.Range("K89").Value = 11214.762
$ws.Range("L89").Value = 416682250
$ws.Range("M89").Value = -5598.762000000001
$ws.Range("N89").Value = -416693482

# Row 105
$ws.Range("H105").Value = 4085.6667
$ws.Range("I105").Value = 3366.6667
$ws.Range("J105").Value = 4485.1113
$ws.Range("K105").Value = 3366.6667
$ws.Range("L105").Value = 4485.1113
$ws.Range("M105").Value = -1619.6667
$ws.Range("N105").Value = -7979.1113

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1717.8182
$ws.Range("I99").Value = 1717.8182
$ws.Range("K99").Value = 1717.8182
$ws.Range("M99").Value = -219.8181999999999

# Row 122
$ws.Range("H122").Value = 1678.1538
$ws.Range("I122").Value = 1734.6666
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5203.9998
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2753.9998
$ws.Range("N122").Value = -7900

# Row 126
$ws.Range("H126").Value = 1717.8182
$ws.Range("I126").Value = 1717.8182
$ws.Range("K126").Value = 5153.4546
$ws.Range("M126").Value = -2683.4546

# Row 141
$ws.Range("H141").Value = 39009.133
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39009.133
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39009.133
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -49369.133

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 19466.5
$ws.Range("J95").Value = 19466.5
$ws.Range("L95").Value = 19466.5
$ws.Range("N95").Value = -24958.5

# Row 102
$ws.Range("H102").Value = 5172.615
$ws.Range("I102").Value = 6154.4
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 6154.4
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = -4532.4
$ws.Range("N102").Value = -5144

# Row 122
$ws.Range("H122").Value = 1954.1904
$ws.Range("I122").Value = 1867
$ws.Range("J122").Value = 2128.5715
$ws.Range("K122").Value = 5601
$ws.Range("L122").Value = 6385.7145
$ws.Range("M122").Value = -3151
$ws.Range("N122").Value = -11285.7145

# Row 126
$ws.Range("H126").Value = 3685.35
$ws.Range("I126").Value = 2234.6667
$ws.Range("J126").Value = 4872.273
$ws.Range("K126").Value = 6704.000100000001
$ws.Range("L126").Value = 14616.819
$ws.Range("M126").Value = -4234.000100000001
$ws.Range("N126").Value = -19556.819

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5340.9062
$ws.Range("I7").Value = 5635.2666
$ws.Range("J7").Value = 5081.1763
$ws.Range("K7").Value = 5635.2666
$ws.Range("L7").Value = 5081.1763
$ws.Range("M7").Value = -5523.2666
$ws.Range("N7").Value = -5305.1763

# Row 122
$ws.Range("H122").Value = 7739
$ws.Range("I122").Value = 6914
$ws.Range("J122").Value = 8357.75
$ws.Range("K122").Value = 20742
$ws.Range("L122").Value = 25073.25
$ws.Range("M122").Value = -18292
$ws.Range("N122").Value = -29973.25

# Row 126
$ws.Range("H126").Value = 5340.9062
$ws.Range("I126").Value = 5635.2666
$ws.Range("J126").Value = 5081.1763
$ws.Range("K126").Value = 16905.7998
$ws.Range("L126").Value = 15243.5289
$ws.Range("M126").Value = -14435.7998
$ws.Range("N126").Value = -20183.5289

# Row 132
$ws.Range("H132").Value = 13897500
$ws.Range("I132").Value = 5481.885
$ws.Range("J132").Value = 50016748
$ws.Range("K132").Value = 16445.655
$ws.Range("L132").Value = 150050244
$ws.Range("M132").Value = -13915.655
$ws.Range("N132").Value = -150055304

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1132.9642
$ws.Range("I136").Value = 1284.4783
$ws.Range("J136").Value = 436
$ws.Range("K136").Value = 3853.4349
$ws.Range("L136").Value = 1308
$ws.Range("M136").Value = -1303.4349
$ws.Range("N136").Value = -6408
